$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Update description of "Greenhouse zone/ irrigation zone" (row 73, column C)
# by clarifying that the zone is not geographically decided.
$ws.Range("C73").Value = "Zone is inside a compartment, so the climate is controlled in the same way but the irrigation. Is the area inside a greenhouse compartment, can be part of all, where the same irrigation and fertigation strategy is applied. Is not geographically decided "

# Add new keyword rows (77-80) describing heating sources / systems as part of
# the roadmap for upcoming ontology versions.
$ws.Range("B77").Value = "Heating source"
$ws.Range("C77").Value = "geothermal energy"

$ws.Range("B78").Value = "Hybrid systems"
$ws.Range("C78").Value = "combines different energy sources for heat. E.g geothermal energy + fuel, or solar collectors"

$ws.Range("B79").Value = "Heat pump"
$ws.Range("C79").Value = "Hybrid system that can be used for cooling and heating "

$ws.Range("B80").Value = "passive system"
$ws.Range("C80").Value = "system where you store energy without any device. Example thick wall to store warmth "

# Adjust the sheet view to match the new scroll/selection position
$win = $excel.ActiveWindow
$win.ScrollRow = 62
$win.ScrollColumn = 1
$ws.Range("C80").Select()
